$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# New row 3 data (the workbook only had a header row + one data row before).
$ws.Range("A3").Value = 112529093
$ws.Range("B3").Value = 89006
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4188
$ws.Range("F3").Value = "Fransig jordstjärna"
$ws.Range("G3").Value = "Geastrum fimbriatum"
$ws.Range("H3").Value = "Fr.:Pers."
# Column I ("Antal") is stored as text "8" in the source file, not a number -
# a leading apostrophe forces Excel to keep it as text.
$ws.Range("I3").Value = "'8"
$ws.Range("J3").Value = "fruktkroppar"
$ws.Range("K3").Value = "teleomorf"
$ws.Range("P3").Value = "Koxutmyren, Upl"
$ws.Range("Q3").Value = 658386
$ws.Range("R3").Value = 6676452
$ws.Range("S3").Value = 50
$ws.Range("T3").Value = "Uppsala"
$ws.Range("U3").Value = "Östhammar"
$ws.Range("V3").Value = "Uppland"
$ws.Range("W3").Value = "Dannemora"
# Dates/times are stored as plain text in this export, not real date values -
# force text with a leading apostrophe so they don't get parsed into serial dates.
$ws.Range("Y3").Value = "'2023-10-05"
$ws.Range("Z3").Value = "09:27"
$ws.Range("AA3").Value = "'2023-10-05"
$ws.Range("AB3").Value = "09:27"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
# AT3/AY3 are present but blank (empty text) in the source row.
$ws.Range("AT3").Value = "'"
$ws.Range("AW3").Value = "Håkan Berglund"
$ws.Range("AX3").Value = "Håkan Berglund"
$ws.Range("AY3").Value = "'"

# Drop the temporary "typed as text" formatting picked up from the leading
# apostrophes above so the new row doesn't carry any explicit cell styling,
# matching the rest of the sheet.
$ws.Rows("3").ClearFormats()
